# Added 4wk low sales check
# Updates the forecast metrics on the "Forecast Comparison" sheet (recomputed
# MyForecast, Inventory Coverage, Stockout Risk, Reorder Urgency and
# Seasonality Index values) and refreshes the dependent roll-up numbers on
# the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---------------------------------------------

# Row 2 (W10)
$ws.Range("D2").Value = 0
$ws.Range("H2").Value = 8.75
$ws.Range("L2").Value = 1.15

# Row 3 (W11)
$ws.Range("D3").Value = 0
$ws.Range("H3").Value = 7.75
$ws.Range("J3").Value = "Normal"
$ws.Range("L3").Value = 0.93

# Row 4 (W12)
$ws.Range("D4").Value = 0
$ws.Range("H4").Value = 6.75
$ws.Range("I4").Value = "Low"
$ws.Range("J4").Value = "Normal"
$ws.Range("L4").Value = 1.03

# Row 5 (W13)
$ws.Range("D5").Value = 1
$ws.Range("H5").Value = 4.6
$ws.Range("I5").Value = "Low"
$ws.Range("J5").Value = "Normal"
$ws.Range("L5").Value = 0.99

# Row 6 (W14)
$ws.Range("D6").Value = 1
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = "Low"
$ws.Range("J6").Value = "Normal"
$ws.Range("L6").Value = 0.98

# Row 7 (W15)
$ws.Range("D7").Value = 1
$ws.Range("H7").Value = 1.71
$ws.Range("I7").Value = "Low"
$ws.Range("J7").Value = "Normal"
$ws.Range("L7").Value = 0.89

# Row 8 (W16)
$ws.Range("D8").Value = 1
$ws.Range("H8").Value = 0.5600000000000001
$ws.Range("I8").Value = "Low"
$ws.Range("L8").Value = 1.03

# Row 9 (W17)
$ws.Range("D9").Value = 2
$ws.Range("L9").Value = 1.06

# Row 10 (W18)
$ws.Range("L10").Value = 0.89

# Row 11 (W19)
$ws.Range("L11").Value = 0.88

# Row 12 (W20)
$ws.Range("L12").Value = 0.87

# Row 13 (W21)
$ws.Range("D13").Value = 3
$ws.Range("L13").Value = 0.97

# Row 14 (W22)
$ws.Range("D14").Value = 3
$ws.Range("L14").Value = 1.11

# Row 15 (W23)
$ws.Range("D15").Value = 3
$ws.Range("L15").Value = 0.97

# Row 16 (W24)
$ws.Range("D16").Value = 4
$ws.Range("L16").Value = 0.85

# Row 17 (W25)
$ws.Range("D17").Value = 4
$ws.Range("L17").Value = 0.88

# --- Summary sheet -----------------------------------------------------------
# Values in column B here are stored as text (not numbers), so force the
# number format to Text first to keep the written constant a text value.

$summary.Range("B9").NumberFormat = "@"
$summary.Range("B9").Value = "39"

$summary.Range("B10").NumberFormat = "@"
$summary.Range("B10").Value = "10"

$summary.Range("B11").NumberFormat = "@"
$summary.Range("B11").Value = "3"

$summary.Range("B12").NumberFormat = "@"
$summary.Range("B12").Value = "5"

$summary.Range("B14").NumberFormat = "@"
$summary.Range("B14").Value = "1"
